# Applies the recorded session of edits to the workbook:
#  1. Select G24 on sheet "1-VLOOKUP" (was A2:E23)
#  2. Select D11 on sheet "7- CHOOSE" (was G15)
#  3. Rename "8- Remove Duplicates " -> "8- REMOVE DUPLICATES"
#  4. Add a new label in E3 on "8- REMOVE DUPLICATES": "Data--> Remove Duplicates "
#  5. Leave "8- REMOVE DUPLICATES" as the active/last-selected sheet

$wb = $excel.ActiveWorkbook

# 1. Update selection on "1-VLOOKUP"
$ws1 = $wb.Worksheets.Item("1-VLOOKUP")
$ws1.Range("G24").Select() | Out-Null

# 2. Update selection on "7- CHOOSE"
$ws7 = $wb.Worksheets.Item("7- CHOOSE")
$ws7.Range("D11").Select() | Out-Null

# 3. Rename the "Remove Duplicates" sheet
$ws8 = $wb.Worksheets.Item("8- Remove Duplicates ")
$ws8.Name = "8- REMOVE DUPLICATES"

# 4. Add the new heading label above the "Remove Duplicates" feature demo
$ws8.Range("E3").Value = "Data--> Remove Duplicates "

# 5. Make this sheet the active tab, as it was the last one touched/saved
$ws8.Activate() | Out-Null
$ws8.Range("I22").Select() | Out-Null
